# Translate the existing Spanish labels/placeholders to English.
# (Writing these in the same order they appear in the target shared-strings
# table keeps the rebuilt sharedStrings.xml ordered exactly like the diff.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "Age"
$ws.Range("A3").Value = "Gender"
$ws.Range("A4").Value = "Height"
$ws.Range("A5").Value = "Weight"

$ws.Range("B1").Value = "{{data.name}}"
$ws.Range("B2").Value = "{{data.age}}"
$ws.Range("B3").Value = "{{data.gender}}"
$ws.Range("B4").Value = "{{data.height}}"
$ws.Range("B5").Value = "{{data.weight}}"

# Add a new row 6 for the "uppercase name" helper, cloning row 5's look
# (font + thin border) so it picks up a style that only differs by font color.
$ws.Range("A5:B5").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)
$ws.Range("A6:B6").Font.Color = 12611584

# Write the template value before the label so the two brand-new shared
# strings land in the order {{uppercase data.name}} (10), Uppercase name (11).
$ws.Range("B6").Value = "{{uppercase data.name}}"
$ws.Range("A6").Value = "Uppercase name"

# Match the workbook's saved selection.
$ws.Range("B4").Select() | Out-Null
